# Weekly data refresh: insert a new Perejil (Vega Central Mapocho de Santiago)
# price-report row at the top of the existing block (row 459), pushing the
# previously-existing rows 459:489 down to 460:490.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 459; Excel shifts rows 459:489 -> 460:490 and
# extends the used range / dimension automatically.
$ws.Rows.Item(459).Insert()

# Populate the newly inserted row with the new weekly record.
$ws.Range("A459").Value = 9
$ws.Range("B459").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C459").Value = "Metropolitana"
$ws.Range("D459").Value = 44931
$ws.Range("E459").Value = 13
$ws.Range("F459").Value = 100112044
$ws.Range("G459").Value = "Perejil"
$ws.Range("H459").Value = "Sin especificar"
$ws.Range("I459").Value = "Primera"
$ws.Range("J459").Value = 70
$ws.Range("K459").Value = 18000
$ws.Range("L459").Value = 18000
$ws.Range("M459").Value = 18000
$ws.Range("N459").Value = "`$/docena de atados"
$ws.Range("O459").Value = "Región Metropolitana"
$ws.Range("P459").Value = 6000
$ws.Range("Q459").Value = 3
$ws.Range("R459").Value = "Hortaliza"
